# Add new status-report rows (65-71) to Sheet1, mirroring the existing
# Date / Hours / Comment table, and update the sheet's scroll/selection
# to reflect the newly-entered data at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Pick up the date-cell formatting (numFmtId 14, style index 4) used by
# the rest of column A by copying the format from the last existing row
# (A64) down onto the new rows, instead of assigning a fresh NumberFormat
# (which would mint a brand-new style entry).
$ws.Range("A64").Copy()
$ws.Range("A65:A71").PasteSpecial(-4122)

# Row 65: 3/23/2010, 2.5h, Group Meeting
$ws.Range("A65").Value = 40260
$ws.Range("B65").Value = 2.5
$ws.Range("C65").Value = "Group Meeting"

# Row 66: 3/23/2010, 1h, Weekly Meeting
$ws.Range("A66").Value = 40260
$ws.Range("B66").Value = 1
$ws.Range("C66").Value = "Weekly Meeting"

# Row 67: 3/24/2010, 3.5h, Physical interface board design
$ws.Range("A67").Value = 40261
$ws.Range("B67").Value = 3.5
$ws.Range("C67").Value = "Physical interface board design"

# Row 68: 3/25/2010, 2h, QNX - Encoder Test
$ws.Range("A68").Value = 40262
$ws.Range("B68").Value = 2
$ws.Range("C68").Value = "QNX - Encoder Test"

# Row 69: 3/26/2010, 0.5h, Skype Meeting
$ws.Range("A69").Value = 40263
$ws.Range("B69").Value = 0.5
$ws.Range("C69").Value = "Skype Meeting"

# Row 70: 3/26/2010, 6h, QNX - Encoder Test
$ws.Range("A70").Value = 40263
$ws.Range("B70").Value = 6
$ws.Range("C70").Value = "QNX - Encoder Test"

# Row 71: 3/27/2010, 8h, QNX - Encoder Test
$ws.Range("A71").Value = 40264
$ws.Range("B71").Value = 8
$ws.Range("C71").Value = "QNX - Encoder Test"

# Scroll the view down to the new bottom of the table and select the
# last entered cell, matching the author's final on-screen position.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 43
$ws.Range("A70").Select()
